# Daily attendance processing - 2025-10-31 08:51:11
#
# The "Recorded By" column (G) lists the users who recorded / touched a
# session, e.g. "System, someone@example.com". This pass normalizes the
# display order of that two-part list by swapping the two names so the
# human recorder is listed first and "System" last (except rows that also
# include the backup@backdoor.com account, which are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$swapped = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text -notlike "*,*") {
        continue
    }

    if ($text -like "*backup@backdoor.com*") {
        continue
    }

    $parts = $text.Split(",")

    if ($parts.Count -ne 2) {
        continue
    }

    $first = $parts[0].Trim()
    $second = $parts[1].Trim()

    $newText = $second + ", " + $first

    $cell.Value = $newText
    $swapped = $swapped + 1
}

Write-Host "Swapped recorder order on" $swapped "rows"
